$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.464.87"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").Value = "2.678.23"
$ws.Range("E3").Value = "  +3.79%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'610.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.50%  "

$ws.Range("D6").Value = "'143.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.40%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.53%  "

$ws.Range("D9").Value = "2.678.01"
$ws.Range("E9").Value = "  +3.79%  "

$ws.Range("E10").Value = "  +0.81%  "

$ws.Range("D11").Value = "'5.61"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "

$ws.Range("E12").Value = "  +0.69%  "

$ws.Range("D13").Value = "'0.358"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.38%  "

$ws.Range("D14").Value = "'27.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.42%  "

$ws.Range("D15").Value = "3.159.31"
$ws.Range("E15").Value = "  +3.80%  "

$ws.Range("D16").Value = "63.344.31"
$ws.Range("E16").Value = "  +0.86%  "

$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").Value = "2.676.31"
$ws.Range("E18").Value = "  +3.90%  "

$ws.Range("E19").Value = "  +3.27%  "

$ws.Range("D20").Value = "'341.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.08%  "

$ws.Range("D21").Value = "'4.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.80%  "

$ws.Range("D22").Value = "'6.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.03%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "'67.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").Value = "'1.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.76%  "

$ws.Range("D26").Value = "'1.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.90%  "

$ws.Range("E27").Value = "  -0.19%  "

$ws.Range("D28").Value = "'8.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.09%  "

$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").Value = "'538.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +17.07%  "

$ws.Range("D31").Value = "'7.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.12%  "

$ws.Range("D32").Value = "'2.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.48%  "

$ws.Range("D33").Value = "'1.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.29%  "

$ws.Range("E34").Value = "  +1.16%  "

$ws.Range("D35").Value = "'172.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.49%  "

$ws.Range("D36").Value = "'5.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.55%  "

$ws.Range("E37").Value = "  +1.71%  "

$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("D39").Value = "'19.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.55%  "

$ws.Range("E40").Value = "  +8.57%  "

$ws.Range("D41").Value = "'176.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.37%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").Value = "'3.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.69%  "

$ws.Range("D44").Value = "'22.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.53%  "

$ws.Range("D45").Value = "'0.0566"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.13%  "

$ws.Range("D46").Value = "'0.635"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.0966"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.45%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0241"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.90%  "

$ws.Range("E49").Value = "  +5.18%  "

$ws.Range("D50").Value = "'1.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.46%  "

$ws.Range("D51").Value = "'11.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.84%  "
